# Update "想去人数" (people interested count) values in the
# 展览 (Exhibitions) and 全部类型 (All Types) sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates
$wsExhibit.Range("F8").Value  = 102
$wsExhibit.Range("F9").Value  = 8648
$wsExhibit.Range("F11").Value = 327
$wsExhibit.Range("F12").Value = 1139
$wsExhibit.Range("F13").Value = 959
$wsExhibit.Range("F14").Value = 100
$wsExhibit.Range("F17").Value = 231
$wsExhibit.Range("F18").Value = 240
$wsExhibit.Range("F19").Value = 63
$wsExhibit.Range("F21").Value = 1000

# 全部类型 sheet updates
$wsAll.Range("F10").Value = 102
$wsAll.Range("F11").Value = 8648
$wsAll.Range("F13").Value = 327
$wsAll.Range("F14").Value = 1139
$wsAll.Range("F15").Value = 959
$wsAll.Range("F16").Value = 100
$wsAll.Range("F19").Value = 231
$wsAll.Range("F20").Value = 240
$wsAll.Range("F21").Value = 63
$wsAll.Range("F23").Value = 1000
